# Update the workbook to reflect the new weekly report run:
#  - refreshed "Report Generated On" timestamp
#  - billed amount / line item counters
#  - all pricing ("Pricing" / H column and TOTAL rows) zeroed out
#  - a new "GND-MD" line item inserted into the Wednesday (08/13/2025) table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header / summary fields -----------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 28

# --- Monday (08/11/2025) pricing -> 0 ----------------------------------
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0

# --- Tuesday (08/12/2025) pricing -> 0 ---------------------------------
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("H35").Value = 0

# --- Wednesday (08/13/2025) pricing -> 0 (rows above the new line item) -
$ws.Range("H40").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("H46").Value = 0

# --- Insert the new "GND-MD" line item as row 47 -----------------------
# This pushes the existing rows 47-51 (INS-15-P-S-C, PIN-XAL-C, POL-45-2,
# SWI-27-CO1-100-H-C, TOTAL) down by one, to rows 48-52. Excel reflows the
# merged cells below automatically (A51:G51->A52:G52, A54:H54->A55:H55,
# A58:G58->A59:G59).
$ws.Rows("47").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# The banded row styling in this table alternates by absolute row number
# (odd rows use the "12/13/14" style triplet, even rows use "9/10/11"),
# so every row from the new one through the TOTAL needs both its content
# and its formatting re-applied to match the new row positions.

# Row 47 (odd) = new GND-MD line item, styled like other odd data rows (e.g. row 45)
$ws.Range("A45:H45").Copy()
$ws.Range("A47:H47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A47").Value = "Point 03"
$ws.Range("B47").Value = "GND-MD"
$ws.Range("C47").Value = "Inst"
$ws.Range("D47").Value = "GND,Wire Mldg Only"
$ws.Range("E47").Value = "EA"
$ws.Range("F47").Value = 2
$ws.Range("H47").Value = 0

# Row 48 (even) = INS-15-P-S-C (was row 47), styled like other even data rows (e.g. row 46)
$ws.Range("A46:H46").Copy()
$ws.Range("A48:H48").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A48").Value = "Point 03"
$ws.Range("B48").Value = "INS-15-P-S-C"
$ws.Range("C48").Value = "Inst"
$ws.Range("D48").Value = "INS,15kV,Pin,Silicon Polymer,Corr"
$ws.Range("E48").Value = "EA"
$ws.Range("F48").Value = 2
$ws.Range("H48").Value = 0

# Row 49 (odd) = PIN-XAL-C (was row 48)
$ws.Range("A45:H45").Copy()
$ws.Range("A49:H49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A49").Value = "Point 03"
$ws.Range("B49").Value = "PIN-XAL-C"
$ws.Range("C49").Value = "Inst"
$ws.Range("D49").Value = "Pin,Crossarm Light,Corrosive"
$ws.Range("E49").Value = "EA"
$ws.Range("F49").Value = 2
$ws.Range("H49").Value = 0

# Row 50 (even) = POL-45-2 (was row 49)
$ws.Range("A46:H46").Copy()
$ws.Range("A50:H50").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A50").Value = "Point 03"
$ws.Range("B50").Value = "POL-45-2"
$ws.Range("C50").Value = "Inst"
$ws.Range("D50").Value = "Pole,45ft,Class 2"
$ws.Range("E50").Value = "EA"
$ws.Range("F50").Value = 1
$ws.Range("H50").Value = 0

# Row 51 (odd) = SWI-27-CO1-100-H-C (was row 50)
$ws.Range("A45:H45").Copy()
$ws.Range("A51:H51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A51").Value = "Point 03"
$ws.Range("B51").Value = "SWI-27-CO1-100-H-C"
$ws.Range("C51").Value = "Inst"
$ws.Range("D51").Value = "SWI,27kV,Line Cutout 1PH,100A,Hook,C"
$ws.Range("E51").Value = "EA"
$ws.Range("F51").Value = 1
$ws.Range("H51").Value = 0

# Row 52 = TOTAL for Wednesday (was row 51); content/format already shifted
# down correctly by the row insert, only the total needs to be zeroed.
$ws.Range("H52").Value = 0

# --- Thursday (08/14/2025) pricing -> 0 (rows shifted down by the insert) -
$ws.Range("H57").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("H59").Value = 0

$wb.Save()
